$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G holds "K" (strikeouts). The save-data generator was changed to
# recompute K from the source game log instead of the old "Strike#" value,
# and std/mean + s_vals were regenerated downstream. Here we just rewrite
# the resulting K values for the existing rows (data rows 2-40).
$gValues = @(
    2, 0, 0, 1, 1, 1, 1, 1, 2, 2,
    2, 1, 1, 0, 1, 1, 1, 1, 1, 3,
    0, 1, 2, 2, 4, 1, 0, 2, 0, 2,
    0, 1, 0, 1, 0, 1, 1, 3, 0
)

for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
}
